$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell updates from the diff. Column D values must remain
# text (as in the source inlineStr cells), so force text format before
# assignment to avoid Excel auto-converting numeric-looking strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.280.41"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.926.46"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "488.46"
$ws.Range("E5").Value = "  +3.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.35"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.734"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("E10").Value = "  +3.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000357"
$ws.Range("E11").Value = "  +5.86%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.67"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.549.38"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.85"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.934.46"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.06"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.387.77"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.09"
$ws.Range("E21").Value = "  +3.33%  "
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.76"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.38"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.51"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.66"
$ws.Range("E25").Value = "  +15.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.07"
$ws.Range("E26").Value = "  +15.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.65"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.94"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.88"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "715.05"
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.131"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0927"
$ws.Range("E34").Value = "  +18.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.68"
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.06"
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.83"
$ws.Range("E37").Value = "  +8.12%  "
$ws.Range("E38").Value = "  -4.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0480"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("E41").Value = "  +15.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.374"
$ws.Range("E42").Value = "  +10.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.09"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("E44").Value = "  +5.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.143"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.42"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "146.13"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0343"
$ws.Range("E50").Value = "  +43.70%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.13"
$ws.Range("E51").Value = "  +0.47%  "
